# Insert a new group of 3 rows (one week of Plátano data for
# Vega Monumental Concepción) above row 320. Excel's row-insert
# shifts all the existing data (rows 320:364) down to (323:367),
# which is exactly what the target workbook needs - the very last
# existing group (D=44273) ends up at rows 365:367 untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A320:T322").EntireRow.Insert()

# Common / constant columns for the new rows.
$commonCols = @{
    1  = 11                             # A Mercado ID
    2  = "Vega Monumental Concepción"   # B Mercado
    3  = "Bíobío"                       # C Región
    5  = 8                              # E Codreg
    6  = "Fruta"                        # F Tipo
    7  = 100108                         # G Producto ID
    8  = "Tropicales y subtropicales"   # H Producto
    9  = 100108006                      # I Categoría ID
    10 = "Plátano"                      # J Categoría
    11 = "Sin especificar"              # K Variedad
    17 = '$/caja 20 kilos'              # Q Unidad de comercialización
    18 = "Ecuador"                      # R Origen
    20 = 20                             # T Kg / unidad
}

# New data: date, quality (L), volume (M), price min/max/avg (N/O/P), $/kg (S)
$newRows = @(
    @{ Row = 320; Fecha = 44505; Calidad = "Maduro";         Volumen = 200; Precio = 13000; PrecioKg = 650 },
    @{ Row = 321; Fecha = 44505; Calidad = "Pintón";          Volumen = 400; Precio = 15000; PrecioKg = 750 },
    @{ Row = 322; Fecha = 44505; Calidad = "Primera Pintón";  Volumen = 400; Precio = 17000; PrecioKg = 850 }
)

foreach ($entry in $newRows) {
    $r = $entry.Row

    foreach ($col in $commonCols.Keys) {
        $ws.Cells.Item($r, $col).Value = $commonCols[$col]
    }

    $ws.Cells.Item($r, 4).Value  = $entry.Fecha      # D Fecha
    $ws.Cells.Item($r, 12).Value = $entry.Calidad    # L Calidad
    $ws.Cells.Item($r, 13).Value = $entry.Volumen    # M Volumen
    $ws.Cells.Item($r, 14).Value = $entry.Precio     # N Precio mínimo
    $ws.Cells.Item($r, 15).Value = $entry.Precio     # O Precio máximo
    $ws.Cells.Item($r, 16).Value = $entry.Precio     # P Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $entry.PrecioKg   # S Precio $/Kg
}
